# Remove the "TestResultExcelFilePath" column (and its values) from the
# "ProcessPayrollForNIMonthly" and "TestReports" worksheets, as the output
# file path should no longer be a fixed/static value (it becomes dynamic).
#
# On both sheets this column is "H" - deleting it shifts every column to the
# right of it one place to the left and Excel automatically renumbers /
# garbage-collects the shared-strings table so the now-unused strings
# ("TestResultExcelFilePath" and the hard-coded output path) disappear.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("ProcessPayrollForNIMonthly")
$ws3.Range("H1").EntireColumn.Delete()

$ws4 = $wb.Worksheets.Item("TestReports")
$ws4.Range("H1").EntireColumn.Delete()

# Reflect the new selection (column H, which now holds what used to be
# column I) on both affected sheets, matching the post-edit workbook state.
$ws3.Activate()
$ws3.Range("H1:H1048576").Select()

$ws4.Activate()
$ws4.Range("H1:H1048576").Select()
